# Updates the cryptos list in sheet1 to reflect the latest scraped values.
# Columns: A=rank(unchanged), B=Coin, C=Link, D=Price, E=Volume(1h)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
    @(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '65.899.96', '  +0.92%  '),
    @(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '3.196.65', '  +0.79%  '),
    @(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '1.00', '  -0.05%  '),
    @(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '599.86', '  +3.67%  '),
    @(6, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '153.71', '  +1.53%  '),
    @(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '1.00', '  -0.03%  '),
    @(8, 'LidoStakedEther', 'https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth', '3.193.75', '  +0.70%  '),
    @(9, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.534', '  +0.65%  '),
    @(10, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.160', '  -1.30%  '),
    @(11, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '6.07', '  -2.32%  '),
    @(12, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.512', '  +1.54%  '),
    @(13, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.0000272', '  +0.82%  '),
    @(14, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '39.08', '  +4.22%  '),
    @(15, 'WrappedliquidstakedEther2.0', 'https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth', '3.719.83', '  +0.79%  '),
    @(16, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '7.45', '  +3.77%  '),
    @(17, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '65.933.12', '  +0.90%  '),
    @(18, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '3.206.38', '  +1.31%  '),
    @(19, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.111', '  -0.02%  '),
    @(20, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '510.33', '  -0.38%  '),
    @(21, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '15.31', '  +2.76%  '),
    @(22, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.740', '  +1.94%  '),
    @(23, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '15.32', '  +0.01%  '),
    @(24, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '8.07', '  +3.14%  '),
    @(25, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '84.84', '  -0.36%  '),
    @(26, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.999', '  -0.21%  '),
    @(27, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '9.34', '  +2.56%  '),
    @(28, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '3.00', '  +2.08%  '),
    @(29, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '2.27', '  +3.49%  '),
    @(30, 'Stacks', 'https://coinranking.com/coin/mMPrMcB7+stacks-stx', '2.88', '  +1.13%  '),
    @(31, 'NEARProtocol', 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near', '6.86', '  +8.71%  '),
    @(32, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '28.11', '  +0.87%  '),
    @(33, 'Mantle', 'https://coinranking.com/coin/BoI4ux0nd+mantle-mnt', '1.22', '  +1.35%  '),
    @(34, 'FirstDigitalUSD', 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd', '1.00', '  +0.03%  '),
    @(35, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '6.55', '  -0.64%  '),
    @(36, 'OKB', 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb', '54.97', '  -0.79%  '),
    @(37, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.0905', '  -0.27%  '),
    @(38, 'Bittensor', 'https://coinranking.com/coin/pgv7xSFi6+bittensor-tao', '485.88', '  +2.24%  '),
    @(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.0418', '  -1.18%  '),
    @(40, 'dogwifhat', 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif', '2.94', '  -4.35%  '),
    @(41, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.89', '  +2.24%  '),
    @(42, 'TheGraph', 'https://coinranking.com/coin/qhd1biQ7M+thegraph-grt', '0.301', '  +5.30%  '),
    @(43, 'Kaspa', 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas', '0.121', '  +2.26%  '),
    @(44, 'PEPE', 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe', '0.0₃0652', '  +8.00%  '),
    @(45, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '2.934.36', '  -4.42%  '),
    @(46, 'Fetch.AI', 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet', '2.43', '  -1.75%  '),
    @(47, 'InjectiveProtocol', 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj', '28.42', '  -2.34%  '),
    @(48, 'USDe', 'https://coinranking.com/coin/exbfr2U-0+usde-usde', '1.00', '  -0.01%  '),
    @(49, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.116', '  +0.93%  '),
    @(50, 'ThetaToken', 'https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta', '2.30', '  +1.84%  '),
    @(51, 'CoreDAO', 'https://coinranking.com/coin/HFvoXUQh4+coredao-core', '2.58', '  +3.96%  ')
)

foreach ($row in $rows) {
    $r = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]

    # Column D ("Price") contains values that often look like plain numbers
    # (e.g. "599.86"); force them to remain text so the stored value matches
    # the scraped string exactly (with its original formatting/precision),
    # then restore the default "Normal" style so no visible formatting change
    # is introduced.
    $dCell = $ws.Cells.Item($r, 4)
    $dCell.NumberFormat = "@"
    $dCell.Value = $row[3]
    $dCell.Style = "Normal"

    $ws.Cells.Item($r, 5).Value = $row[4]
}
